$p = $ppt.ActivePresentation

# --- Slide 1: subtitle "Jesse" " " "Rosenthal" -> single run "Jesse Rosenthal" ---
# The paragraph is: <a:br/><a:br/><a:r>Jesse</a:r><a:r> </a:r><a:r>Rosenthal</a:r>
# Grow the first text run in place (preserves <a:pPr> and the run's <a:rPr/>),
# then clear out the now-redundant trailing runs.
$s1 = $p.Slides.Item(1)
$subtitleTr = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitleRun1 = $subtitleTr.Characters(3, 5)      # "Jesse"
$subtitleRun1.Text = "Jesse Rosenthal"
$subtitleRun2 = $subtitleTr.Characters(18, 1)     # " "
$subtitleRun2.Text = ""
$subtitleRun3 = $subtitleTr.Characters(18, 9)     # "Rosenthal"
$subtitleRun3.Text = ""

# --- Slide 1 speaker notes: "Some" " " "speaker" " " "notes" -> "Some speaker notes" ---
$notesTr = $s1.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesTr.Text = "Some speaker notes"

# --- Slide 2: title "A" " " "header" -> single run "A header" ---
$s2 = $p.Slides.Item(2)
$titleTr = $s2.Shapes.Item(1).TextFrame.TextRange
$titleRun1 = $titleTr.Characters(1, 1)            # "A"
$titleRun1.Text = "A header"
$titleRun2 = $titleTr.Characters(9, 1)            # " "
$titleRun2.Text = ""
$titleRun3 = $titleTr.Characters(9, 6)            # "header"
$titleRun3.Text = ""
